$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 7 ("gru_cls_sep" / 0.871) merges into row 6 by deleting old row 6
# ("linear" / 0.9412) and shifting row 7 up.
$ws.Rows.Item(6).Delete() | Out-Null

# New "Optimizations" column.
$ws.Range("D1").Value = "Optimizations"
$ws.Range("D2:D6").Value = "da: drop_col - dk: None - summarize: False"

# Match the header formatting used on A1:C1 (bold font, centered/top
# alignment) but with only left/right borders on the new column.
$ws.Range("A1").Copy() | Out-Null
$ws.Range("D1").PasteSpecial(-4122) | Out-Null
$ws.Range("D1").Borders.Item(9).LineStyle = -4142
$ws.Range("D1").Borders.Item(8).LineStyle = -4142

$ws.Columns.Item(4).ColumnWidth = 38.5

$ws.Range("D6").Select() | Out-Null
